$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Both columns end up with the same stored width (15.42578125 chars in the
# OOXML <col> width attribute). The ColumnWidth COM property is quantized
# to whole pixels (Normal-style font metrics), so 14.6 "characters" is the
# input that lands on the nearest reachable stored width (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 14.6

$values = @(
    @(-0.36880911882937539, 0.36801746440559668),
    @(-0.1838642692842356, 0.18260907490364531),
    @(-0.08919606939359781, 0.088814063588834813),
    @(-0.15080422509542046, 0.15007845669844144),
    @(-0.14407845793370466, 0.14262895462986247),
    @(-0.04172534779295578, 0.041699713743543754),
    @(-0.021699715215046211, 0.021640790914158714),
    @(-0.0016407923898746901, 0.0015380292963564912),
    @(0.0044619694242378216, -0.0045783710485940432),
    @(0.010578369770009033, -0.010593153433376301),
    @(0.0019765998859320177, -0.0019751960174225758),
    @(0.0079751947387132205, -0.0079831172581350529),
    @(0.013983115982158623, -0.013991443741268483),
    @(0.025991442382482077, -0.026025738738669979),
    @(0.032025737468191373, -0.032088182623798644),
    @(0.038088181359484441, -0.038239998874767256),
    @(0.044239997624659466, -0.044351738860677337),
    @(-0.1137154746440423, 0.11359006077330847),
    @(-0.027095659276784811, 0.027012451224365019),
    @(-0.018012452470147267, 0.018004084218645744),
    @(-0.0090040854661346259, 0.0089999987515083646),
    @(-0.0844130067645672, 0.084205784951661755),
    @(-0.075205786216142378, 0.074860104399731675),
    @(-0.042124518446779113, 0.041999998242409298),
    @(-0.094903608108822368, 0.09466613940400137),
    @(-0.088666140680793148, 0.088356183559543666),
    @(-0.082356184843045632, 0.081280693773574608),
    @(-0.075280695084983584, 0.074535888037061682),
    @(-0.062535889451837079, 0.062168570639036247),
    @(-0.042168572175321994, 0.042017449617713254),
    @(-0.027017451096202905, 0.027000279000887417),
    @(-0.0060002805658063707, 0.0059999986464260857)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
